$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen), P (Precio $/Kg).
$updates = @(
    @{ Row = 2; D = 44230; J = 250; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 3; D = 44208; J = 160; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 4; D = 44236; J = 180; K = 4000; L = 4500; M = 4167; O = 'Región Metropolitana'; P = 260 }
    @{ Row = 5; D = 44875; J = 90; K = 7000; L = 7000; M = 7000; O = 'Provincia de Quillota'; P = 438 }
    @{ Row = 6; D = 44882; J = 70; K = 7000; L = 7000; M = 7000; O = 'Provincia de Quillota'; P = 438 }
    @{ Row = 7; D = 44210; J = 340; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 8; D = 44232; J = 250; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 9; D = 44873; J = 250; K = 8000; L = 8000; M = 8000; O = 'Provincia de Quillota'; P = 500 }
    @{ Row = 10; D = 44188; J = 210; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 11; D = 44204; J = 430; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 12; D = 44855; J = 70; K = 6000; L = 7000; M = 6500; O = 'Provincia de Quillota'; P = 406 }
    @{ Row = 13; D = 44251; J = 120; K = 5000; L = 5000; M = 5000; O = 'Región Metropolitana'; P = 312 }
    @{ Row = 14; D = 44215; J = 250; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 15; D = 44186; J = 160; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 16; D = 44883; J = 180; K = 7000; L = 8000; M = 7500; O = 'Provincia de Quillota'; P = 469 }
    @{ Row = 17; D = 44189; J = 250; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 18; D = 44187; J = 160; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 19; D = 44231; J = 250; K = 5000; L = 6000; M = 5500; O = 'Provincia de Quillota'; P = 344 }
    @{ Row = 20; D = 44292; J = 90; K = 6000; L = 6000; M = 6000; O = 'Región Metropolitana'; P = 375 }
    @{ Row = 21; D = 44846; J = 250; K = 5000; L = 5000; M = 5000; O = 'Provincia de Quillota'; P = 312 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 10).Value = $u.J
    $ws.Cells.Item($u.Row, 11).Value = $u.K
    $ws.Cells.Item($u.Row, 12).Value = $u.L
    $ws.Cells.Item($u.Row, 13).Value = $u.M
    $ws.Cells.Item($u.Row, 15).Value = $u.O
    $ws.Cells.Item($u.Row, 16).Value = $u.P
}
